# Updates cryptos list price/volume figures (GitHub Actions scheduled refresh),
# plus the insertion of "BabyDogeCoin" which bumps Algorand/Cronos down one row
# and drops the former last row (EnergySwap) off the bottom of the table.
# Values that look numeric (e.g. "215.47") are written with a leading "'" so
# Excel keeps them as literal text instead of re-parsing them as numbers,
# exactly like the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.875.30"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.664.32"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'215.47"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'0.534"
$ws.Range("E6").Value = "  +5.44%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "'0.0620"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("E10").Value = "  +3.55%  "
$ws.Range("D11").Value = "'0.0898"
$ws.Range("E11").Value = "  +3.88%  "
$ws.Range("D12").Value = "1.898.73"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "1.661.54"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "'66.14"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "26.877.19"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "'232.94"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").Value = "'7.96"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "'4.40"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("D24").Value = "'9.15"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").Value = "'146.00"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "'0.115"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("E32").Value = "  +2.06%  "
$ws.Range("D33").Value = "1.459.81"
$ws.Range("E33").Value = "  -4.75%  "
$ws.Range("E34").Value = "  +3.80%  "
$ws.Range("E35").Value = "  +3.47%  "
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").Value = "'0.900"
$ws.Range("E38").Value = "  +2.12%  "
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("E40").Value = "  -3.69%  "
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("E42").Value = "  -2.91%  "
$ws.Range("E43").Value = "  +6.02%  "
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").Value = "1.808.45"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").Value = "'0.780"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").Value = "'90.40"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0104"
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.102"
$ws.Range("E50").Value = "  +4.75%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0506"
$ws.Range("E51").Value = "  +0.40%  "
